$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 5 with a sample value in column A, mirroring the
# "Sample row added" entry already present in row 4.
$ws.Range("A5").Value = "Mysample"

# Move the active selection to A6 (the row following the newly added row)
$ws.Range("A6").Select()
